$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.078.38'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.817.54'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5900'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2727'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06794'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07524'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '1.826.48'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.634'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6216'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009482'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.29%  '
$ws.Range('D17').Value = '28.844.81'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.420'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.65%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.88%  '
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.754'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.007'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '154.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.782'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1263'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06339'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.93%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.409'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.429'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.706'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.678'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.686'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.050'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.18%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6353'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.534'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.741'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01703'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.363'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('D40').Value = '1.133.44'
$ws.Range('E40').Value = '  -8.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8623'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.24%  '
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.966.92'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.569'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.96%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4538'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05483'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.229'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.006'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.24%  '
